# Auto-generated edit script: applies Word auto-proofing-style run splits
# (w:proofErr spellStart/spellEnd/gramStart/gramEnd markers) plus the
# 'Задорожний Костянтин' content addition, matching the target diff.

$d = $word.ActiveDocument

# Row1 col1 para2: 'Test Suite Description' -> word-by-word proofErr
$r = $d.Content
$ok = $r.Find.Execute('Test Suite Description', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok) { throw 'Find failed for: Test Suite Description' }
$r.Expand(4) | Out-Null
$r.InsertXML('<w:p w14:paraId="178295CD" w14:textId="77777777" w:rsidR="00F13932" w:rsidRPr="00725406" w:rsidRDefault="00F13932" w:rsidP="0036395B" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:right="-214"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Test</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Suite</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Description</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

# Row2 col1 para2: 'Name of ' + 'Modules/ prj / ModulesZadorozhny' -> split w/ proofErr
$r = $d.Content
$ok = $r.Find.Execute('Name of ', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok) { throw 'Find failed for: Name of ' }
$r.Expand(4) | Out-Null
$r.InsertXML('<w:p w14:paraId="092C3CD9" w14:textId="6E5BA727" w:rsidR="00F13932" w:rsidRPr="00F13932" w:rsidRDefault="00F13932" w:rsidP="0036395B" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:right="-214"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>of</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Modules/ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>prj</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> / </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>ModulesZadorozhny</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

# Row2 col2 para: 's_calculation()' -> split w/ proofErr/gramErr
$r = $d.Content
$ok = $r.Find.Execute('s_calculation()', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok) { throw 'Find failed for: s_calculation()' }
$r.Expand(4) | Out-Null
$r.InsertXML('<w:p w14:paraId="01270570" w14:textId="679F038A" w:rsidR="00F13932" w:rsidRPr="00B93355" w:rsidRDefault="00B073A9" w:rsidP="0036395B" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:left="104"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:sz w:val="25"/><w:szCs w:val="25"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00B93355"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:sz w:val="25"/><w:szCs w:val="25"/><w:lang w:val="ru-RU"/></w:rPr><w:t>s_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:sz w:val="25"/><w:szCs w:val="25"/><w:lang w:val="ru-RU"/></w:rPr><w:t>calculation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:sz w:val="25"/><w:szCs w:val="25"/><w:lang w:val="ru-RU"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:sz w:val="25"/><w:szCs w:val="25"/><w:lang w:val="ru-RU"/></w:rPr><w:t>)</w:t></w:r><w:r w:rsidRPr="00B93355"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:sz w:val="25"/><w:szCs w:val="25"/><w:lang w:val="en-US"/></w:rPr><w:t>;</w:t></w:r></w:p>')

# Row3 col1 para2: 'Level of Testing' -> word-by-word proofErr
$r = $d.Content
$ok = $r.Find.Execute('Level of Testing', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok) { throw 'Find failed for: Level of Testing' }
$r.Expand(4) | Out-Null
$r.InsertXML('<w:p w14:paraId="49E3C7A6" w14:textId="77777777" w:rsidR="00F13932" w:rsidRPr="00725406" w:rsidRDefault="00F13932" w:rsidP="0036395B" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:right="-214"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Level</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>of</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Testing</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

# Row3 col2 para: last run ' Testing' -> ' ' + 'Testing' w/ proofErr
$r = $d.Content
$ok = $r.Find.Execute('Unit Testing', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok) { throw 'Find failed for: Unit Testing' }
$r.Expand(4) | Out-Null
$r.InsertXML('<w:p w14:paraId="426F06F1" w14:textId="0975E00F" w:rsidR="00F13932" w:rsidRPr="00725406" w:rsidRDefault="00F13932" w:rsidP="0036395B" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:left="104" w:right="-218"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:sz w:val="25"/><w:szCs w:val="25"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t>модульний</w:t></w:r><w:r w:rsidRPr="00725406"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:sz w:val="25"/><w:szCs w:val="25"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">  /  </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:lang w:val="en-US"/></w:rPr><w:t>Unit</w:t></w:r><w:r w:rsidRPr="00725406"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/></w:rPr><w:t>Testing</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

# Row4 col1 para1: 'Автор тест-сьюта ' -> split w/ proofErr
$r = $d.Content
$ok = $r.Find.Execute('Автор тест-сьюта ', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok) { throw 'Find failed for: Автор тест-сьюта ' }
$r.Expand(4) | Out-Null
$r.InsertXML('<w:p w14:paraId="5C00D54D" w14:textId="77777777" w:rsidR="00F13932" w:rsidRPr="00725406" w:rsidRDefault="00F13932" w:rsidP="0036395B" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:right="-214"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="00725406"><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Автор тест-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>сьюта</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>')

# Row4 col1 para2: 'Test Suite Author' -> word-by-word proofErr
$r = $d.Content
$ok = $r.Find.Execute('Test Suite Author', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok) { throw 'Find failed for: Test Suite Author' }
$r.Expand(4) | Out-Null
$r.InsertXML('<w:p w14:paraId="2703C231" w14:textId="77777777" w:rsidR="00F13932" w:rsidRPr="00725406" w:rsidRDefault="00F13932" w:rsidP="0036395B" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:right="-214"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Test</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Suite</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Author</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

# Row5 col1 para2: 'Implementer' -> wrap w/ proofErr
$r = $d.Content
$ok = $r.Find.Execute('Implementer', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok) { throw 'Find failed for: Implementer' }
$r.Expand(4) | Out-Null
$r.InsertXML('<w:p w14:paraId="189A3A99" w14:textId="77777777" w:rsidR="00F13932" w:rsidRPr="00725406" w:rsidRDefault="00F13932" w:rsidP="0036395B" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:right="-214"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00725406"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Implementer</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

# Row5 col2 para: empty -> add 'Задорожний Костянтин' runs
$t1 = $d.Tables.Item(1)
$targetCell = $t1.Cell(5,2)
$cr = $targetCell.Range.Paragraphs.Item(1).Range
$cr.InsertXML('<w:p w14:paraId="6FC6C3D7" w14:textId="6E4BD8B4" w:rsidR="00F13932" w:rsidRPr="00725406" w:rsidRDefault="00F13932" w:rsidP="0036395B" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:left="104"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t>З</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t>адорожний Костянтин</w:t></w:r></w:p>')

# Second table, row2 col2 para2: 'X = 120 ;' -> split w/ gramErr
$r = $d.Content
$ok = $r.Find.Execute('X = 120 ;', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok) { throw 'Find failed for: X = 120 ;' }
$r.Expand(4) | Out-Null
$r.InsertXML('<w:p w14:paraId="7522CEB0" w14:textId="77777777" w:rsidR="00F13932" w:rsidRDefault="00F13932" w:rsidP="0036395B" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">X = </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>120 ;</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

Write-Output "done"